# Update the "Förändrad" (Changed) date column (C2:C6) from 45212 to 45221
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
